$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 23.75
$ws.Range("E2").Value = 26.04999923706055
$ws.Range("F2").Value = 26.89999961853028
$ws.Range("G2").Value = 22.60000038146973
$ws.Range("H2").Value = 168459019
$ws.Range("I2").Value = "OKTA"

$ws.Range("D3").Value = 23.75
$ws.Range("E3").Value = 26.04999923706055
$ws.Range("F3").Value = 26.89999961853028
$ws.Range("G3").Value = 22.60000038146973
$ws.Range("H3").Value = 168459019
$ws.Range("I3").Value = "OKTA"

$ws.Range("D4").Value = 23.75
$ws.Range("E4").Value = 26.04999923706055
$ws.Range("F4").Value = 26.89999961853028
$ws.Range("G4").Value = 22.60000038146973
$ws.Range("H4").Value = 168459019
$ws.Range("I4").Value = "OKTA"

$ws.Range("D5").Value = 23.75
$ws.Range("E5").Value = 26.04999923706055
$ws.Range("F5").Value = 26.89999961853028
$ws.Range("G5").Value = 22.60000038146973
$ws.Range("H5").Value = 168459019
$ws.Range("I5").Value = "OKTA"

$ws.Range("D6").Value = 23.75
$ws.Range("E6").Value = 26.04999923706055
$ws.Range("F6").Value = 26.89999961853028
$ws.Range("G6").Value = 22.60000038146973
$ws.Range("H6").Value = 168459019
$ws.Range("I6").Value = "OKTA"

$ws.Range("D7").Value = 26.5
$ws.Range("E7").Value = 26.10000038146973
$ws.Range("F7").Value = 26.79999923706055
$ws.Range("G7").Value = 22.67000007629395
$ws.Range("H7").Value = 168459019
$ws.Range("I7").Value = "OKTA"

$ws.Range("D8").Value = 22.28000068664551
$ws.Range("E8").Value = 26.98999977111816
$ws.Range("F8").Value = 27.47999954223633
$ws.Range("G8").Value = 21.52000045776367
$ws.Range("H8").Value = 168459019
$ws.Range("I8").Value = "OKTA"

$ws.Range("D9").Value = 29
$ws.Range("E9").Value = 29.20999908447266
$ws.Range("F9").Value = 31.79999923706055
$ws.Range("G9").Value = 27.96999931335449
$ws.Range("H9").Value = 168459019
$ws.Range("I9").Value = "OKTA"

$ws.Range("D10").Value = 29.38999938964844
$ws.Range("E10").Value = 38.59000015258789
$ws.Range("F10").Value = 39.13000106811523
$ws.Range("G10").Value = 27.70999908447266
$ws.Range("H10").Value = 168459019
$ws.Range("I10").Value = "OKTA"

$ws.Range("D11").Value = 42.79999923706055
$ws.Range("E11").Value = 56.20999908447266
$ws.Range("F11").Value = 56.91899871826172
$ws.Range("G11").Value = 42.63999938964844
$ws.Range("H11").Value = 168459019
$ws.Range("I11").Value = "OKTA"

$ws.Range("D12").Value = 50
$ws.Range("E12").Value = 61.83000183105469
$ws.Range("F12").Value = 62.38999938964844
$ws.Range("G12").Value = 50
$ws.Range("H12").Value = 168459019
$ws.Range("I12").Value = "OKTA"

$ws.Range("D13").Value = 59.2400016784668
$ws.Range("E13").Value = 63.65000152587891
$ws.Range("F13").Value = 64
$ws.Range("G13").Value = 41.88000106811523
$ws.Range("H13").Value = 168459019
$ws.Range("I13").Value = "OKTA"

$ws.Range("D14").Value = 82.20999908447266
$ws.Range("E14").Value = 84.87999725341797
$ws.Range("F14").Value = 87.71900177001953
$ws.Range("G14").Value = 78.76000213623047
$ws.Range("H14").Value = 168459019
$ws.Range("I14").Value = "OKTA"

$ws.Range("D15").Value = 105
$ws.Range("E15").Value = 113.2200012207031
$ws.Range("F15").Value = 119.9599990844727
$ws.Range("G15").Value = 98.91000366210938
$ws.Range("H15").Value = 168459019
$ws.Range("I15").Value = "OKTA"

$ws.Range("D16").Value = 132
$ws.Range("E16").Value = 126.5
$ws.Range("F16").Value = 138.9389953613281
$ws.Range("G16").Value = 121.3399963378906
$ws.Range("H16").Value = 168459019
$ws.Range("I16").Value = "OKTA"

$ws.Range("D17").Value = 109.5599975585938
$ws.Range("E17").Value = 129.7799987792969
$ws.Range("F17").Value = 131.5500030517578
$ws.Range("G17").Value = 104.3300018310547
$ws.Range("H17").Value = 168459019
$ws.Range("I17").Value = "OKTA"

$ws.Range("D18").Value = 128.75
$ws.Range("E18").Value = 128.0599975585938
$ws.Range("F18").Value = 142.9799957275391
$ws.Range("G18").Value = 118.5800018310547
$ws.Range("H18").Value = 168459019
$ws.Range("I18").Value = "OKTA"

$ws.Range("D19").Value = 148.6900024414062
$ws.Range("E19").Value = 195.5800018310547
$ws.Range("F19").Value = 196.8999938964844
$ws.Range("G19").Value = 147.2400054931641
$ws.Range("H19").Value = 168459019
$ws.Range("I19").Value = "OKTA"

$ws.Range("D20").Value = 221.6900024414062
$ws.Range("E20").Value = 215.3699951171875
$ws.Range("F20").Value = 226.8899993896484
$ws.Range("G20").Value = 193.7100067138672
$ws.Range("H20").Value = 168459019
$ws.Range("I20").Value = "OKTA"

$ws.Range("D21").Value = 210.0050048828125
$ws.Range("E21").Value = 245.0399932861328
$ws.Range("F21").Value = 246.6399993896484
$ws.Range("G21").Value = 200.6199951171875
$ws.Range("H21").Value = 168459019
$ws.Range("I21").Value = "OKTA"

$ws.Range("D22").Value = 260.7099914550781
$ws.Range("E22").Value = 261.4500122070312
$ws.Range("F22").Value = 294
$ws.Range("G22").Value = 252
$ws.Range("H22").Value = 168459019
$ws.Range("I22").Value = "OKTA"

$ws.Range("D23").Value = 269.9500122070312
$ws.Range("E23").Value = 222.4400024414062
$ws.Range("F23").Value = 270.4599914550781
$ws.Range("G23").Value = 215.6199951171875
$ws.Range("H23").Value = 168459019
$ws.Range("I23").Value = "OKTA"

$ws.Range("D24").Value = 248
$ws.Range("E24").Value = 263.6000061035156
$ws.Range("F24").Value = 266.9469909667969
$ws.Range("G24").Value = 226.0599975585937
$ws.Range("H24").Value = 168459019
$ws.Range("I24").Value = "OKTA"

$ws.Range("D25").Value = 247.1999969482422
$ws.Range("E25").Value = 215.229995727539
$ws.Range("F25").Value = 272.2699890136719
$ws.Range("G25").Value = 208.2400054931641
$ws.Range("H25").Value = 168459019
$ws.Range("I25").Value = "OKTA"

$ws.Range("D26").Value = 199.3999938964844
$ws.Range("E26").Value = 182.8399963378907
$ws.Range("F26").Value = 203.7899932861328
$ws.Range("G26").Value = 152.5099945068359
$ws.Range("H26").Value = 168459019
$ws.Range("I26").Value = "OKTA"

$ws.Range("D27").Value = 118.4100036621094
$ws.Range("E27").Value = 83.05000305175781
$ws.Range("F27").Value = 125.3000030517578
$ws.Range("G27").Value = 77.01000213623047
$ws.Range("H27").Value = 168459019
$ws.Range("I27").Value = "OKTA"

$ws.Range("D28").Value = 96.91000366210938
$ws.Range("E28").Value = 91.40000152587891
$ws.Range("F28").Value = 110.9400024414062
$ws.Range("G28").Value = 88.22000122070312
$ws.Range("H28").Value = 168459019
$ws.Range("I28").Value = "OKTA"

$ws.Range("D29").Value = 58.11999893188477
$ws.Range("E29").Value = 53.31999969482422
$ws.Range("F29").Value = 59.13000106811523
$ws.Range("G29").Value = 44.11999893188477
$ws.Range("H29").Value = 168459019
$ws.Range("I29").Value = "OKTA"

$ws.Range("D30").Value = 73.5
$ws.Range("E30").Value = 71.29000091552734
$ws.Range("F30").Value = 82.09999847412109
$ws.Range("G30").Value = 70.12000274658203
$ws.Range("H30").Value = 168459019
$ws.Range("I30").Value = "OKTA"

$ws.Range("D31").Value = 67.41000366210938
$ws.Range("E31").Value = 90.90000152587891
$ws.Range("F31").Value = 91.5
$ws.Range("G31").Value = 66.81999969482422
$ws.Range("H31").Value = 168459019
$ws.Range("I31").Value = "OKTA"

$ws.Range("D32").Value = 76.37000274658203
$ws.Range("E32").Value = 83.51000213623047
$ws.Range("F32").Value = 86.5
$ws.Range("G32").Value = 68.79000091552734
$ws.Range("H32").Value = 168459019
$ws.Range("I32").Value = "OKTA"

$ws.Range("D33").Value = 67.77999877929688
$ws.Range("E33").Value = 67.05000305175781
$ws.Range("F33").Value = 73.29000091552734
$ws.Range("G33").Value = 65.04000091552734
$ws.Range("H33").Value = 168459019
$ws.Range("I33").Value = "OKTA"

$ws.Range("D34").Value = 83.83999633789062
$ws.Range("E34").Value = 107.3000030517578
$ws.Range("F34").Value = 112.0800018310547
$ws.Range("G34").Value = 79.34999847412109
$ws.Range("H34").Value = 168459019
$ws.Range("I34").Value = "OKTA"

$ws.Range("D35").Value = 94.16000366210938
$ws.Range("E35").Value = 88.68000030517578
$ws.Range("F35").Value = 104.109001159668
$ws.Range("G35").Value = 86.25399780273438
$ws.Range("H35").Value = 168459019
$ws.Range("I35").Value = "OKTA"

$ws.Range("D36").Value = 94.31999969482422
$ws.Range("E36").Value = 78.73000335693359
$ws.Range("F36").Value = 99.91000366210938
$ws.Range("G36").Value = 77.63999938964844
$ws.Range("H36").Value = 168459019
$ws.Range("I36").Value = "OKTA"

$ws.Range("D37").Value = 72.02999877929688
$ws.Range("E37").Value = 77.55999755859375
$ws.Range("F37").Value = 80.70999908447266
$ws.Range("G37").Value = 70.91999816894531
$ws.Range("H37").Value = 168459019
$ws.Range("I37").Value = "OKTA"

$ws.Range("D38").Value = 92.5
$ws.Range("E38").Value = 90.48999786376952
$ws.Range("F38").Value = 100.5999984741211
$ws.Range("G38").Value = 87.37999725341797
$ws.Range("H38").Value = 168459019
$ws.Range("I38").Value = "OKTA"

$ws.Range("D39").Value = 113.3399963378906
$ws.Range("E39").Value = 103.1699981689453
$ws.Range("F39").Value = 127.5670013427734
$ws.Range("G39").Value = 102.1500015258789
$ws.Range("H39").Value = 168459019
$ws.Range("I39").Value = "OKTA"

$ws.Range("D40").Value = 96.41000366210938
$ws.Range("E40").Value = 92.76999664306641
$ws.Range("F40").Value = 98.87000274658205
$ws.Range("G40").Value = 87.23999786376953
$ws.Range("H40").Value = 168459019
$ws.Range("I40").Value = "OKTA"

